# Apply the "new database input file standards" edit to the VRES profiles
# workbook. Both sheets (ScenarioA, ScenarioB) share an identical header
# layout (rows 1-7), so the same changes are applied to each worksheet:
#
#   - H4: long description text -> "Capacity" (column header code label)
#   - E5: (empty) -> "Technology of the capacity factor" (description)
#   - F5: "Which package this node belongs to"
#         -> "Which data package this belongs to" (wording fix)
#   - H5: (empty) -> "Capacityfactor for each VRES of this technology at
#         this node " (description), formatted like B5:G5 but without
#         word-wrap (a new, very similar cell style)
#   - H7:AE7: "[p.u.]" -> "[%, 0-1]" (unit row, capacity factors are a
#         0-1 ratio, not a per-unit voltage/power quantity)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- H4: rename the long leftover description to the short "Capacity" label
    $ws.Range("H4").Value2 = "Capacity"

    # --- E5: add the missing description for the Technology column
    $ws.Range("E5").Value2 = "Technology of the capacity factor"

    # --- F5: reword the Data Package description
    $ws.Range("F5").Value2 = "Which data package this belongs to"

    # --- H5: add the description for the capacity-factor data block.
    # Copy the formatting of the neighbouring description cell (B5) so we
    # get the same italic font / fill / left+top alignment, then turn off
    # word-wrap (H5 uses a dedicated un-wrapped variant of that style).
    $ws.Range("B5").Copy() | Out-Null
    $ws.Range("H5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("H5").WrapText = $false
    $ws.Range("H5").Value2 = "Capacityfactor for each VRES of this technology at this node "

    # --- H7:AE7: update the unit label for every capacity-factor column
    $ws.Range("H7:AE7").Value2 = "[%, 0-1]"

}

$excel.CutCopyMode = 0
